$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current values A13:A41 into an array
$values = @()
for ($r = 13; $r -le 41; $r++) {
    $values += $ws.Cells.Item($r, 1).Value2
}

# Rotate up by one: value from row N moves to row N-1, first value wraps to last row (A41)
$first = $values[0]
for ($i = 0; $i -lt ($values.Length - 1); $i++) {
    $ws.Cells.Item(13 + $i, 1).Value = $values[$i + 1]
}
$ws.Cells.Item(41, 1).Value = $first
